# Applies the "Lab 3" edit: adds rolling STDEV.S(sample-of-5) formulas in
# columns F/G beneath every existing rolling-AVERAGE row, updates the
# worksheet zoom/selection, and leaves the dimension to be recomputed from
# the newly-populated cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each block of 45 data rows (5 reps x 9 days) already carries an
# AVERAGE(D..:D..) / AVERAGE(E..:E..) pair in columns F/G on the last row
# of every 5-row rep. Immediately below every one of those rows we add the
# matching sample standard deviation over the same 5-row window. The final
# block's last STDEV pair lands one row past the previous last used row,
# so it becomes a brand new row 199.
$blockStarts = 10, 58, 106, 154
$blockEnds   = 54, 102, 150, 198

for ($b = 0; $b -lt $blockStarts.Length; $b++) {
    $start = $blockStarts[$b]
    $end   = $blockEnds[$b]

    for ($row = $start + 5; $row -le $end + 1; $row += 5) {
        $winFirst = $row - 5
        $winLast  = $row - 1

        $ws.Range("F$row").Formula = "=STDEV.S(D${winFirst}:D${winLast})"
        $ws.Range("G$row").Formula = "=STDEV.S(E${winFirst}:E${winLast})"
    }
}

# Zoom out and move the active selection, matching the saved view state.
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("I76").Select()
